$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Cells.Item(19, 8).Value = 4111.9473
$ws_ALC.Cells.Item(19, 9).Value = 2173.5557
$ws_ALC.Cells.Item(19, 10).Value = 5856.5
$ws_ALC.Cells.Item(19, 11).Value = 2173.5557
$ws_ALC.Cells.Item(19, 12).Value = 5856.5
$ws_ALC.Cells.Item(19, 13).Value = -1998.5557
$ws_ALC.Cells.Item(19, 14).Value = -6206.5
$ws_ALC.Cells.Item(40, 8).Value = 5492.2856
$ws_ALC.Cells.Item(40, 10).Value = 9284.571
$ws_ALC.Cells.Item(40, 12).Value = 9284.571
$ws_ALC.Cells.Item(40, 14).Value = -9634.571
$ws_ALC.Cells.Item(62, 8).Value = 1200
$ws_ALC.Cells.Item(62, 10).Value = 1200
$ws_ALC.Cells.Item(62, 12).Value = 1200
$ws_ALC.Cells.Item(62, 14).Value = -2448
$ws_ALC.Cells.Item(65, 8).Value = 1200
$ws_ALC.Cells.Item(65, 10).Value = 1200
$ws_ALC.Cells.Item(65, 12).Value = 6000
$ws_ALC.Cells.Item(65, 14).Value = -12240
$ws_ALC.Cells.Item(76, 8).Value = 2750
$ws_ALC.Cells.Item(76, 9).Value = 3000
$ws_ALC.Cells.Item(76, 11).Value = 3000
$ws_ALC.Cells.Item(76, 13).Value = -2685
$ws_ALC.Cells.Item(79, 8).Value = 2750
$ws_ALC.Cells.Item(79, 9).Value = 3000
$ws_ALC.Cells.Item(79, 11).Value = 3000
$ws_ALC.Cells.Item(79, 13).Value = -1908
$ws_ALC.Cells.Item(98, 8).Value = 1565.8334
$ws_ALC.Cells.Item(98, 9).Value = 1565.8334
$ws_ALC.Cells.Item(98, 11).Value = 1565.8334
$ws_ALC.Cells.Item(98, 13).Value = -67.83339999999998
$ws_ALC.Cells.Item(106, 8).Value = 5777.6
$ws_ALC.Cells.Item(106, 9).Value = 5777.6
$ws_ALC.Cells.Item(106, 11).Value = 5777.6
$ws_ALC.Cells.Item(106, 13).Value = -5146.6
$ws_ALC.Cells.Item(107, 8).Value = 901.3611
$ws_ALC.Cells.Item(107, 9).Value = 905.4516
$ws_ALC.Cells.Item(107, 10).Value = 876
$ws_ALC.Cells.Item(107, 11).Value = 905.4516
$ws_ALC.Cells.Item(107, 12).Value = 876
$ws_ALC.Cells.Item(107, 13).Value = 1014.5484
$ws_ALC.Cells.Item(107, 14).Value = -4716
$ws_ALC.Cells.Item(111, 8).Value = 2793
$ws_ALC.Cells.Item(111, 9).Value = 2791.25
$ws_ALC.Cells.Item(111, 11).Value = 8373.75
$ws_ALC.Cells.Item(111, 13).Value = -5306.75
$ws_ALC.Cells.Item(122, 8).Value = 1565.8334
$ws_ALC.Cells.Item(122, 9).Value = 1565.8334
$ws_ALC.Cells.Item(122, 11).Value = 4697.5002
$ws_ALC.Cells.Item(122, 13).Value = -2247.5002
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Cells.Item(109, 8).Value = 66160.625
$ws_ARM.Cells.Item(109, 10).Value = 67055
$ws_ARM.Cells.Item(109, 12).Value = 67055
$ws_ARM.Cells.Item(109, 14).Value = -69829
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Cells.Item(86, 8).Value = 200001500
$ws_BSM.Cells.Item(86, 9).Value = 500001600
$ws_BSM.Cells.Item(86, 11).Value = 500001600
$ws_BSM.Cells.Item(86, 13).Value = -500000477
$ws_BSM.Cells.Item(89, 8).Value = 200001500
$ws_BSM.Cells.Item(89, 9).Value = 500001600
$ws_BSM.Cells.Item(89, 11).Value = 2500008000
$ws_BSM.Cells.Item(89, 13).Value = -2500002384
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Cells.Item(19, 8).Value = 474.375
$ws_CRP.Cells.Item(19, 9).Value = 209
$ws_CRP.Cells.Item(19, 11).Value = 209
$ws_CRP.Cells.Item(19, 13).Value = -39
$ws_CRP.Cells.Item(24, 8).Value = 474.375
$ws_CRP.Cells.Item(24, 9).Value = 209
$ws_CRP.Cells.Item(24, 11).Value = 209
$ws_CRP.Cells.Item(24, 13).Value = -39
$ws_CRP.Cells.Item(74, 8).Value = 96219.336
$ws_CRP.Cells.Item(74, 10).Value = 96219.336
$ws_CRP.Cells.Item(74, 12).Value = 96219.336
$ws_CRP.Cells.Item(74, 14).Value = -97967.336
$ws_CRP.Cells.Item(77, 8).Value = 96219.336
$ws_CRP.Cells.Item(77, 10).Value = 96219.336
$ws_CRP.Cells.Item(77, 12).Value = 288658.008
$ws_CRP.Cells.Item(77, 14).Value = -297394.008
$ws_CRP.Cells.Item(107, 8).Value = 1442.8572
$ws_CRP.Cells.Item(107, 9).Value = 1267.1052
$ws_CRP.Cells.Item(107, 10).Value = 1813.8889
$ws_CRP.Cells.Item(107, 11).Value = 1267.1052
$ws_CRP.Cells.Item(107, 12).Value = 1813.8889
$ws_CRP.Cells.Item(107, 13).Value = 652.8948
$ws_CRP.Cells.Item(107, 14).Value = -5653.8889
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Cells.Item(4, 8).Value = 208571
$ws_CUL.Cells.Item(4, 9).Value = 259073.02
$ws_CUL.Cells.Item(4, 11).Value = 777219.0599999999
$ws_CUL.Cells.Item(4, 13).Value = -777107.0599999999
$ws_CUL.Cells.Item(6, 8).Value = 113.63636
$ws_CUL.Cells.Item(6, 9).Value = 71.77778000000001
$ws_CUL.Cells.Item(6, 10).Value = 302
$ws_CUL.Cells.Item(6, 11).Value = 215.33334
$ws_CUL.Cells.Item(6, 12).Value = 906
$ws_CUL.Cells.Item(6, 13).Value = -102.33334
$ws_CUL.Cells.Item(6, 14).Value = -1132
$ws_CUL.Cells.Item(100, 8).Value = 6999.5
$ws_CUL.Cells.Item(100, 9).Value = 6999.5
$ws_CUL.Cells.Item(100, 11).Value = 20998.5
$ws_CUL.Cells.Item(100, 13).Value = -20187.5
$ws_CUL.Cells.Item(140, 8).Value = 1532.4762
$ws_CUL.Cells.Item(140, 9).Value = 1121.2222
$ws_CUL.Cells.Item(140, 11).Value = 3363.6666
$ws_CUL.Cells.Item(140, 13).Value = 1816.3334
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Cells.Item(9, 8).Value = 1870.8
$ws_GSM.Cells.Item(9, 10).Value = 2627.5
$ws_GSM.Cells.Item(9, 12).Value = 2627.5
$ws_GSM.Cells.Item(9, 14).Value = -2967.5
$ws_GSM.Cells.Item(80, 8).Value = 2648.5
$ws_GSM.Cells.Item(80, 9).Value = 2049.5
$ws_GSM.Cells.Item(80, 10).Value = 3247.5
$ws_GSM.Cells.Item(80, 11).Value = 2049.5
$ws_GSM.Cells.Item(80, 12).Value = 3247.5
$ws_GSM.Cells.Item(80, 13).Value = -1051.5
$ws_GSM.Cells.Item(80, 14).Value = -5243.5
$ws_GSM.Cells.Item(83, 8).Value = 2648.5
$ws_GSM.Cells.Item(83, 9).Value = 2049.5
$ws_GSM.Cells.Item(83, 10).Value = 3247.5
$ws_GSM.Cells.Item(83, 11).Value = 10247.5
$ws_GSM.Cells.Item(83, 12).Value = 16237.5
$ws_GSM.Cells.Item(83, 13).Value = -5255.5
$ws_GSM.Cells.Item(83, 14).Value = -26221.5
$ws_GSM.Cells.Item(97, 8).Value = 2267.8
$ws_GSM.Cells.Item(97, 9).Value = 2313.889
$ws_GSM.Cells.Item(97, 11).Value = 2313.889
$ws_GSM.Cells.Item(97, 13).Value = -1817.889
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Cells.Item(2, 8).Value = 150008500
$ws_LTW.Cells.Item(2, 10).Value = 12750
$ws_LTW.Cells.Item(2, 12).Value = 12750
$ws_LTW.Cells.Item(2, 14).Value = -12974
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Cells.Item(21, 8).Value = 30017
$ws_WVR.Cells.Item(21, 10).Value = 30017
$ws_WVR.Cells.Item(21, 12).Value = 30017
$ws_WVR.Cells.Item(21, 14).Value = -30487
$ws_WVR.Cells.Item(35, 8).Value = 30017
$ws_WVR.Cells.Item(35, 10).Value = 30017
$ws_WVR.Cells.Item(35, 12).Value = 30017
$ws_WVR.Cells.Item(35, 14).Value = -30597
$ws_WVR.Cells.Item(81, 8).Value = 3248308.2
$ws_WVR.Cells.Item(81, 9).Value = 5103172
$ws_WVR.Cells.Item(81, 10).Value = 2296.5
$ws_WVR.Cells.Item(81, 11).Value = 10206344
$ws_WVR.Cells.Item(81, 12).Value = 4593
$ws_WVR.Cells.Item(81, 13).Value = -10205283
$ws_WVR.Cells.Item(81, 14).Value = -6715
$ws_WVR.Cells.Item(84, 8).Value = 3248308.2
$ws_WVR.Cells.Item(84, 9).Value = 5103172
$ws_WVR.Cells.Item(84, 10).Value = 2296.5
$ws_WVR.Cells.Item(84, 11).Value = 51031720
$ws_WVR.Cells.Item(84, 12).Value = 22965
$ws_WVR.Cells.Item(84, 13).Value = -51026416
$ws_WVR.Cells.Item(84, 14).Value = -33573
$ws_WVR.Cells.Item(122, 8).Value = 3490.9333
$ws_WVR.Cells.Item(122, 9).Value = 3383.1428
$ws_WVR.Cells.Item(122, 11).Value = 10149.4284
$ws_WVR.Cells.Item(122, 13).Value = -7699.428400000001
$ws_WVR.Cells.Item(132, 8).Value = 1516.8611
$ws_WVR.Cells.Item(132, 9).Value = 1377.4839
$ws_WVR.Cells.Item(132, 10).Value = 2381
$ws_WVR.Cells.Item(132, 11).Value = 4132.4517
$ws_WVR.Cells.Item(132, 12).Value = 7143
$ws_WVR.Cells.Item(132, 13).Value = -1602.4517
$ws_WVR.Cells.Item(132, 14).Value = -12203
$ws_WVR.Cells.Item(136, 8).Value = 1255.8148
$ws_WVR.Cells.Item(136, 9).Value = 1086.5
$ws_WVR.Cells.Item(136, 10).Value = 2000.8
$ws_WVR.Cells.Item(136, 11).Value = 3259.5
$ws_WVR.Cells.Item(136, 12).Value = 6002.4
$ws_WVR.Cells.Item(136, 13).Value = -709.5
$ws_WVR.Cells.Item(136, 14).Value = -11102.4
